# Minor flow cleanup and commenting
# Applies small numeric corrections to the proposed_rhna_allocation sheet
# (column C = "existing need" and column F adjustments), as part of a
# minor data cleanup pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("proposed_rhna_allocation")

# Column C corrections
$ws.Range("C8").Value   = 618
$ws.Range("C20").Value  = 745
$ws.Range("C23").Value  = 197
$ws.Range("C31").Value  = 325
$ws.Range("C43").Value  = 19
$ws.Range("C52").Value  = 59
$ws.Range("C64").Value  = 104
$ws.Range("C71").Value  = 406
$ws.Range("C86").Value  = 1598
$ws.Range("C88").Value  = 8689
$ws.Range("C95").Value  = 448
$ws.Range("C110").Value = 221

# Column F corrections
$ws.Range("F23").Value  = 343
$ws.Range("F31").Value  = 566
$ws.Range("F41").Value  = 261
$ws.Range("F48").Value  = 1100
$ws.Range("F52").Value  = 102
$ws.Range("F64").Value  = 181
$ws.Range("F74").Value  = 871
$ws.Range("F110").Value = 385
